$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: update timestamp, high_node, num_high_nodes, P_per_hi, Q_per_hi ---
$ws.Cells.Item(4, 2).Value = 44519.68430940972
$ws.Cells.Item(4, 5).Value = "Bus_151"
$ws.Cells.Item(4, 7).Value = 11
$ws.Cells.Item(4, 8).Value = -707.89484178924647
$ws.Cells.Item(4, 9).Value = -342.62110342599527

# --- Row 5: update timestamp, type, high_node, num_high_nodes, P/Q_per_hi, P/Q_per_lo; C5 becomes numeric 123 ---
$ws.Cells.Item(5, 2).Value = 44519.690525034719
$ws.Cells.Item(5, 3).Value = 123
$ws.Cells.Item(5, 4).Value = "Undervoltage"
$ws.Cells.Item(5, 5).Value = "Bus_151"
$ws.Cells.Item(5, 7).Value = 11
$ws.Cells.Item(5, 8).Value = -707.89484178924567
$ws.Cells.Item(5, 9).Value = -342.62110342599487
$ws.Cells.Item(5, 11).Value = 657.93028129657807
$ws.Cells.Item(5, 12).Value = 318.43825614754383

# --- New sample run rows 6-9 ---
# Clone the existing row formatting (col A bordered/bold/centered style,
# col B timestamp number format) from rows 2:5 down onto the 4 new rows,
# so we reuse the existing styles instead of synthesizing new ones.
$ws.Range("A2:B5").Copy()
$ws.Range("A6:B9").PasteSpecial(-4122)
$excel.CutCopyMode = $false | Out-Null

$newRows = @(
    @{ Row=6; A=4; B=44519.687113819447;  C=123; D="Overvoltage";  E="Bus_151"; F="Bus_79"; G=11; H=-707.89484178924647;  I=-342.62110342599527;  J=13; K=260.0672937048343;  L=125.8725701531398 },
    @{ Row=7; A=5; B=44519.688869687503;  C=123; D="Undervoltage"; E="Bus_151"; F="Bus_79"; G=11; H=-707.89484178924567;  I=-342.62110342599487;  J=13; K=698.07536731297603; L=337.86847777948037 },
    @{ Row=8; A=6; B=44519.691084571758;  C=123; D="Overvoltage";  E="Bus_48";  F="Bus_79"; G=8;  H=-1279.2671069477101;  I=-619.16527976269151;  J=13; K=260.0672937048343;  L=125.8725701531398 },
    @{ Row=9; A=7; B=44519.692170556627;  C=123; D="Undervoltage"; E="Bus_48";  F="Bus_79"; G=8;  H=-1279.267106947708;   I=-619.16527976269072;  J=13; K=876.0161472162838;  L=423.99181525268131 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
}

# --- Update sheet selection to match the saved cursor position ---
$ws.Range("G5").Select() | Out-Null
